$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 19 new rows for the additional DEC_0771..DEC_0789 test-data rows
# (pushes the old "footer" rows 56-60 down to 75-79).
$ws.Rows("56:74").Insert()

# Insert 2 more blank separator rows before the old blank (style 4) row,
# matching the author's extra spacing (old 58-60 end up at 79-81).
$ws.Rows("77:78").Insert()

# Fill the 19 new data rows (56-74) with the same pattern used by the
# existing rows above them (e.g. row 55): TC number, password, plantilla,
# and SIN_DATO placeholders, with the script id incrementing in column A.
for ($i = 0; $i -lt 19; $i++) {
    $row = 56 + $i
    $num = 771 + $i
    $dec = "DEC_0$num"

    $ws.Cells.Item($row, 1).Value = $dec
    $ws.Cells.Item($row, 2).Value = "18092588-0"
    # Leading apostrophe preserves the "number stored as text" (quote
    # prefix) formatting that column C already used for this value.
    $ws.Cells.Item($row, 3).Value = "'sebA`$1357"

    for ($col = 4; $col -le 10; $col++) {
        $ws.Cells.Item($row, $col).Value = "SIN_DATO"
    }
}

# Mirror the final cursor position left by the editing session.
$ws.Range("B74").Select()
$excel.ActiveWindow.ScrollRow = 61
